$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = "U2FsdGVkX18CeviQ24hIEERMC5lllAkCyTq8qB+K1J8ELLtppxyKsUFxOE+nntq9/WW38YXJQgmprt6+xWNAmDhQSQC6lExckmwooXEjTNnJ3+9TvopXVeVD4S0/efSymJkoILmDy9RhiqXoSqODsFf1mw71a6OLivpxueg/q+qqYdvq9yTI5iNQnz8y5ZvghwX8Tnm229QpdNMizeUnveiSMqY7iXdMtRUQvaeCOZiQmAcsaONQgiZHeeIJfpeKgFkOVDkheWibIj5j8MfxlLI/fFNKAEvF/SJnCDktvpPdSOUz0PHB2+E7GdBSsjK5"

$ws.Range("B3").Value = $newValue
$ws.Range("B5").Value = $newValue
$ws.Range("B7").Value = $newValue
